$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row content: B1 and C1 change text, D1 is cleared
$ws.Range("B1").Value = "點餐者"
$ws.Range("C1").Value = "時間"
$ws.Range("D1").ClearContents()

# Set column widths
$ws.Columns.Item(2).ColumnWidth = 14.75
$ws.Columns.Item(3).ColumnWidth = 21.5

# Default row height
$ws.Cells.Item(1, 1).RowHeight = 15.75

# Update selection
$ws.Range("F3").Select()
